$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-19 down to 4-20
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new market data entry
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44921
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = "Frutos de pepita"
$ws.Range("I3").Value = 100104004
$ws.Range("J3").Value = "Níspero"
$ws.Range("K3").Value = "Californiana(o)"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 7 kilos"
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 2143
$ws.Range("T3").Value = 7
